$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fueltrade")

# Update the output commodity columns (Comm/Comm1/Comm2, columns D:F) from H2 to H2GC
# for every trade-link row, then rename the Tech column (G) TB_H2_* identifiers to
# TB_H2GC_* to match. The commodity cells are updated first so the new "H2GC" shared
# string is allocated before the old TB_H2_* strings are freed/renamed.
for ($r = 4; $r -le 11; $r++) {
    $ws.Cells.Item($r, 4).Value = "H2GC"
    $ws.Cells.Item($r, 5).Value = "H2GC"
    $ws.Cells.Item($r, 6).Value = "H2GC"
}

$techNames = @(
    "TB_H2GC_DKISLBH_DKE_01",
    "TB_H2GC_DKISLBH_DKE_02",
    "TB_H2GC_DKISL1_DKW_01",
    "TB_H2GC_DKISL1_DKW_02",
    "TB_H2GC_DKISL2_DKW_01",
    "TB_H2GC_DKISL2_DKW_02",
    "TB_H2GC_DKISL3_DKW_01",
    "TB_H2GC_DKISL3_DKW_02"
)
for ($i = 0; $i -lt $techNames.Length; $i++) {
    $r = 4 + $i
    $ws.Cells.Item($r, 7).Value = $techNames[$i]
}

# Make "Fueltrade" the active/selected tab (it was "BI" before).
$ws.Activate()
